$p = $ppt.ActivePresentation
Write-Output "HasTitleMaster: $($p.HasTitleMaster)"
Write-Output "HasHandoutMaster: $($p.HasHandoutMaster)"
$hm = $p.HandoutMaster
if ($hm -eq $null) { Write-Output "handout null" } else { Write-Output "handout ok" }
